$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "29.309.18"
Set-TextValue $ws "E2" "  +0.60%  "
Set-TextValue $ws "D3" "1.874.32"
Set-TextValue $ws "E3" "  +0.82%  "
Set-TextValue $ws "D4" "0.9996"
Set-TextValue $ws "E4" "  -0.16%  "
Set-TextValue $ws "D5" "0.7126"
Set-TextValue $ws "E5" "  +0.04%  "
Set-TextValue $ws "D6" "242.79"
Set-TextValue $ws "E6" "  +1.21%  "
Set-TextValue $ws "D7" "0.9998"
Set-TextValue $ws "E7" "  -0.10%  "
Set-TextValue $ws "D8" "0.3113"
Set-TextValue $ws "E8" "  +1.60%  "
Set-TextValue $ws "D9" "0.07756"
Set-TextValue $ws "E9" "  +0.71%  "
Set-TextValue $ws "D10" "25.12"
Set-TextValue $ws "E10" "  +1.17%  "
Set-TextValue $ws "D11" "0.08476"
Set-TextValue $ws "E11" "  +2.97%  "
Set-TextValue $ws "D12" "1.865.54"
Set-TextValue $ws "E12" "  +0.45%  "
Set-TextValue $ws "D13" "5.214"
Set-TextValue $ws "E13" "  +0.15%  "
Set-TextValue $ws "D14" "0.7121"
Set-TextValue $ws "E14" "  -0.29%  "
Set-TextValue $ws "D15" "91.37"
Set-TextValue $ws "E15" "  +1.35%  "
Set-TextValue $ws "B16" "ShibaInu"
Set-TextValue $ws "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D16" "0.000008364"
Set-TextValue $ws "E16" "  +7.58%  "
Set-TextValue $ws "B17" "WrappedBTC"
Set-TextValue $ws "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws "D17" "29.302.19"
Set-TextValue $ws "E17" "  +0.59%  "
Set-TextValue $ws "D18" "5.989"
Set-TextValue $ws "E18" "  +2.57%  "
Set-TextValue $ws "D19" "243.01"
Set-TextValue $ws "E19" "  -0.04%  "
Set-TextValue $ws "B20" "Avalanche"
Set-TextValue $ws "C20" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws "D20" "13.23"
Set-TextValue $ws "E20" "  +1.01%  "
Set-TextValue $ws "B21" "WrappedliquidstakedEther2.0"
Set-TextValue $ws "C21" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws "D21" "2.125.40"
Set-TextValue $ws "E21" "  +1.28%  "
Set-TextValue $ws "D22" "0.9996"
Set-TextValue $ws "E22" "  -0.11%  "
Set-TextValue $ws "D23" "7.801"
Set-TextValue $ws "E23" "  -1.40%  "
Set-TextValue $ws "D24" "1.0000"
Set-TextValue $ws "E24" "  -0.19%  "
Set-TextValue $ws "D25" "0.1623"
Set-TextValue $ws "E25" "  +2.77%  "
Set-TextValue $ws "D26" "162.94"
Set-TextValue $ws "E26" "  +0.68%  "
Set-TextValue $ws "D27" "9.025"
Set-TextValue $ws "E27" "  +1.61%  "
Set-TextValue $ws "D28" "18.50"
Set-TextValue $ws "D29" "1.509"
Set-TextValue $ws "E29" "  +1.12%  "
Set-TextValue $ws "D30" "4.421"
Set-TextValue $ws "E30" "  +2.10%  "
Set-TextValue $ws "D31" "4.336"
Set-TextValue $ws "E31" "  +6.39%  "
Set-TextValue $ws "D32" "1.275"
Set-TextValue $ws "E32" "  -3.03%  "
Set-TextValue $ws "E33" "  +1.78%  "
Set-TextValue $ws "D34" "1.922"
Set-TextValue $ws "E34" "  +0.86%  "
Set-TextValue $ws "D35" "1.174"
Set-TextValue $ws "E35" "  +0.33%  "
Set-TextValue $ws "D36" "0.7433"
Set-TextValue $ws "E36" "  +2.26%  "
Set-TextValue $ws "D37" "2.683"
Set-TextValue $ws "E37" "  +0.22%  "
Set-TextValue $ws "D38" "0.01861"
Set-TextValue $ws "E38" "  +1.05%  "
Set-TextValue $ws "D39" "2.720"
Set-TextValue $ws "E39" "  +1.04%  "
Set-TextValue $ws "D40" "1.165.15"
Set-TextValue $ws "E40" "  +1.40%  "
Set-TextValue $ws "D41" "6.348"
Set-TextValue $ws "E41" "  +4.64%  "
Set-TextValue $ws "D42" "0.8895"
Set-TextValue $ws "E42" "  -0.89%  "
Set-TextValue $ws "D43" "72.94"
Set-TextValue $ws "E43" "  +1.32%  "
Set-TextValue $ws "D44" "106.55"
Set-TextValue $ws "E44" "  +5.03%  "
Set-TextValue $ws "E45" "  -0.05%  "
Set-TextValue $ws "D46" "2.020.81"
Set-TextValue $ws "E46" "  +0.70%  "
Set-TextValue $ws "D47" "1.807"
Set-TextValue $ws "E47" "  +2.85%  "
Set-TextValue $ws "D48" "0.5193"
Set-TextValue $ws "E48" "  -1.57%  "
Set-TextValue $ws "D50" "9.379"
Set-TextValue $ws "E50" "  +1.42%  "
Set-TextValue $ws "D51" "0.4304"
Set-TextValue $ws "E51" "  +1.98%  "
